$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B,C,D values for rows 2-9 (re-recorded data for 2a)
$ws.Range("B2").Value = -0.003550
$ws.Range("C2").Value = 5.772366
$ws.Range("D2").Value = 0.000000

$ws.Range("B3").Value = -0.007835
$ws.Range("C3").Value = 5.723294
$ws.Range("D3").Value = 0.049072

$ws.Range("B4").Value = -0.012120
$ws.Range("C4").Value = 5.772366
$ws.Range("D4").Value = 0.000000

$ws.Range("B5").Value = -0.016405
$ws.Range("C5").Value = 5.772366
$ws.Range("D5").Value = 0.000000

$ws.Range("B6").Value = -0.020690
$ws.Range("C6").Value = 5.772366
$ws.Range("D6").Value = 0.000000

$ws.Range("B7").Value = -0.024975
$ws.Range("C7").Value = 5.772366
$ws.Range("D7").Value = 0.000000

$ws.Range("B8").Value = -0.029260
$ws.Range("C8").Value = 5.968656
$ws.Range("D8").Value = 0.196290

$ws.Range("B9").Value = -0.033545
$ws.Range("C9").Value = 6.606598
$ws.Range("D9").Value = 0.834232

# Remove rows 10-12, which no longer exist in the re-recorded data
$ws.Range("A10:D12").Delete()
